$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 5
$ws.Range("D11").Value = 5
$ws.Range("E12").Value = 5
$ws.Range("D14").Value = 5
$ws.Range("C15").Value = 5
$ws.Range("E15").Value = 5
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 5
$ws.Range("E18").Value = 5
$ws.Range("C21").Value = 5
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = 5
$ws.Range("D24").Value = 5
$ws.Range("C28").Value = 5

$ws.Range("E22").Select()
